$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "219.45", "4.440") are preserved exactly instead of being coerced
# into floating point numbers (which would corrupt formatting / trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.275.13'
$ws.Cells.Item(2, 5).Value = '  +0.73%  '
$ws.Cells.Item(3, 4).Value = '1.657.81'
$ws.Cells.Item(3, 5).Value = '  +0.26%  '
$ws.Cells.Item(4, 5).Value = '  +0.57%  '
$ws.Cells.Item(5, 4).Value = '219.45'
$ws.Cells.Item(5, 5).Value = '  +2.12%  '
$ws.Cells.Item(6, 4).Value = '0.5215'
$ws.Cells.Item(6, 5).Value = '  -0.46%  '
$ws.Cells.Item(7, 4).Value = '1.006'
$ws.Cells.Item(7, 5).Value = '  +0.54%  '
$ws.Cells.Item(8, 5).Value = '  +0.94%  '
$ws.Cells.Item(9, 4).Value = '0.06326'
$ws.Cells.Item(9, 5).Value = '  -0.80%  '
$ws.Cells.Item(10, 4).Value = '21.38'
$ws.Cells.Item(10, 5).Value = '  +2.66%  '
$ws.Cells.Item(11, 4).Value = '0.07753'
$ws.Cells.Item(11, 5).Value = '  +0.22%  '
$ws.Cells.Item(12, 4).Value = '1.669.58'
$ws.Cells.Item(12, 5).Value = '  +1.58%  '
$ws.Cells.Item(13, 4).Value = '4.440'
$ws.Cells.Item(13, 5).Value = '  -0.38%  '
$ws.Cells.Item(14, 4).Value = '0.5493'
$ws.Cells.Item(14, 5).Value = '  -0.42%  '
$ws.Cells.Item(15, 5).Value = '  -0.56%  '
$ws.Cells.Item(16, 4).Value = '64.99'
$ws.Cells.Item(16, 5).Value = '  -0.05%  '
$ws.Cells.Item(17, 4).Value = '26.284.84'
$ws.Cells.Item(17, 5).Value = '  +0.77%  '
$ws.Cells.Item(18, 4).Value = '1.007'
$ws.Cells.Item(18, 5).Value = '  +0.58%  '
$ws.Cells.Item(19, 4).Value = '4.705'
$ws.Cells.Item(19, 5).Value = '  -0.97%  '
$ws.Cells.Item(20, 4).Value = '191.71'
$ws.Cells.Item(20, 5).Value = '  +0.78%  '
$ws.Cells.Item(21, 4).Value = '10.21'
$ws.Cells.Item(21, 5).Value = '  -0.42%  '
$ws.Cells.Item(22, 4).Value = '6.218'
$ws.Cells.Item(22, 5).Value = '  -2.16%  '
$ws.Cells.Item(23, 4).Value = '1.008'
$ws.Cells.Item(23, 5).Value = '  +0.71%  '
$ws.Cells.Item(24, 4).Value = '138.95'
$ws.Cells.Item(24, 5).Value = '  -2.93%  '
$ws.Cells.Item(25, 4).Value = '0.1254'
$ws.Cells.Item(25, 5).Value = '  +0.40%  '
$ws.Cells.Item(26, 4).Value = '7.311'
$ws.Cells.Item(26, 5).Value = '  -1.36%  '
$ws.Cells.Item(27, 4).Value = '16.08'
$ws.Cells.Item(27, 5).Value = '  +0.56%  '
$ws.Cells.Item(28, 4).Value = '1.425'
$ws.Cells.Item(28, 5).Value = '  +0.29%  '
$ws.Cells.Item(29, 4).Value = '0.06063'
$ws.Cells.Item(29, 5).Value = '  +2.07%  '
$ws.Cells.Item(30, 4).Value = '1.289'
$ws.Cells.Item(30, 5).Value = '  +2.12%  '
$ws.Cells.Item(31, 4).Value = '3.559'
$ws.Cells.Item(31, 5).Value = '  +3.45%  '
$ws.Cells.Item(32, 4).Value = '3.379'
$ws.Cells.Item(32, 5).Value = '  -1.12%  '
$ws.Cells.Item(33, 4).Value = '1.664'
$ws.Cells.Item(33, 5).Value = '  +0.68%  '
$ws.Cells.Item(34, 4).Value = '0.9912'
$ws.Cells.Item(34, 5).Value = '  -0.57%  '
$ws.Cells.Item(35, 4).Value = '2.428'
$ws.Cells.Item(36, 4).Value = '2.775'
$ws.Cells.Item(36, 5).Value = '  +0.67%  '
$ws.Cells.Item(37, 4).Value = '0.5982'
$ws.Cells.Item(37, 5).Value = '  +6.09%  '
$ws.Cells.Item(38, 4).Value = '0.01599'
$ws.Cells.Item(38, 5).Value = '  -0.29%  '
$ws.Cells.Item(39, 4).Value = '5.989'
$ws.Cells.Item(39, 5).Value = '  +1.89%  '
$ws.Cells.Item(40, 4).Value = '1.074.04'
$ws.Cells.Item(40, 5).Value = '  +4.60%  '
$ws.Cells.Item(41, 4).Value = '0.8529'
$ws.Cells.Item(41, 5).Value = '  -0.57%  '
$ws.Cells.Item(42, 5).Value = '  +0.41%  '
$ws.Cells.Item(43, 4).Value = '99.86'
$ws.Cells.Item(43, 5).Value = '  +0.51%  '
$ws.Cells.Item(44, 4).Value = '1.802.29'
$ws.Cells.Item(44, 5).Value = '  +0.14%  '
$ws.Cells.Item(45, 4).Value = '57.66'
$ws.Cells.Item(45, 5).Value = '  +3.25%  '
$ws.Cells.Item(46, 5).Value = '  +0.27%  '
$ws.Cells.Item(47, 4).Value = '1.004'
$ws.Cells.Item(48, 4).Value = '8.077'
$ws.Cells.Item(48, 5).Value = '  +0.12%  '
$ws.Cells.Item(49, 4).Value = '0.05193'
$ws.Cells.Item(49, 5).Value = '  +0.84%  '
$ws.Cells.Item(50, 4).Value = '1.473'
$ws.Cells.Item(50, 5).Value = '  +6.00%  '
$ws.Cells.Item(51, 4).Value = '0.4233'
$ws.Cells.Item(51, 5).Value = '  +0.47%  '
